$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3. This shifts the existing row 3 (and its
# formatting) down to row 4, so the previous week's reading is preserved
# there, and leaves row 3 free for the new week's reading.
$ws.Rows.Item(3).Insert()

# Row 3: new weekly reading (same descriptive fields as surrounding rows,
# with updated date and price/volume values)
$ws.Range("A3").Value = 8
$ws.Range("B3").Value = "Terminal La Palmera de La Serena"
$ws.Range("C3").Value = "Coquimbo"
$ws.Range("D3").Value = 45237
$ws.Range("E3").Value = 4
$ws.Range("F3").Value = 300000000
$ws.Range("G3").Value = "Espárragos"
$ws.Range("H3").Value = "Sin especificar"
$ws.Range("I3").Value = "Primera"
$ws.Range("J3").Value = 800
$ws.Range("K3").Value = 1900
$ws.Range("L3").Value = 2000
$ws.Range("M3").Value = 1950
$ws.Range("N3").Value = "`$/kilo"
$ws.Range("O3").Value = "Provincia de Linares"
$ws.Range("P3").Value = 1950
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = "Hortaliza"

$ws.Range("D3").NumberFormat = $ws.Range("D4").NumberFormat

$wb.Save()
